# This edit permutes the content of data rows 2..22 (row 1 is the header).
# Each destination row ends up with the exact former content of another
# source row (a derangement / permutation over rows 2..22), columns A..AY.
#
# Mapping: destRow -> sourceRow (i.e. destRow receives what sourceRow used to contain)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$perm = @{2=16; 3=21; 4=20; 5=22; 6=10; 7=11; 8=17; 9=6; 10=2; 11=19; 12=18; 13=14; 14=7; 15=12; 16=3; 17=5; 18=8; 19=15; 20=9; 21=13; 22=4}

$firstCol = 1
$lastCol = 51

# 1) Snapshot every source row's values first, since destinations overlap
#    with sources (it's a full permutation) and writes must not clobber
#    data that is still needed for a later read.
$snapshot = @{}
for ($r = 2; $r -le 22; $r++) {
    $rowVals = @()
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $rowVals += $ws.Cells.Item($r, $col).Value()
    }
    $snapshot[$r] = $rowVals
}

# 2) Write each destination row from the snapshot of its mapped source row.
#    Columns Y and AA hold plain "yyyy-mm-dd" text in this sheet (no real
#    date typing anywhere in the workbook). Assigning such a string straight
#    to .Value lets Excel's type-sniffing reinterpret it as a real date
#    serial, so force those two columns to Text first to keep them as the
#    original literal strings.
$dateTextCols = @(25, 27)
for ($destRow = 2; $destRow -le 22; $destRow++) {
    $srcRow = $perm[$destRow]
    $rowVals = $snapshot[$srcRow]
    foreach ($dc in $dateTextCols) {
        $ws.Cells.Item($destRow, $dc).NumberFormat = "@"
    }
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $ws.Cells.Item($destRow, $col).Value = $rowVals[$col - $firstCol]
    }
    foreach ($dc in $dateTextCols) {
        $ws.Cells.Item($destRow, $dc).NumberFormat = "General"
    }
}
